$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (Ambiente row) with the new environment's data ---
$ws.Range("A2").Value() = "i-preproducciongestion.segurossura.com.ar"

# Replace the B2 hyperlink (remove old one, set new display text + target)
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.Delete()
    }
}
$ws.Range("B2").Value() = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do")
# Re-adding a hyperlink re-applies the hyperlink cell style as a brand new
# style entry; reuse the existing hyperlink style (shared with B3/B4) instead.
$ws.Range("B2").Style = $ws.Range("B3").Style

$ws.Range("C2").Value() = "su"
# Leading apostrophe keeps this a text value (preserving the leading zero
# and the existing quote-prefix cell style) instead of becoming a number.
$ws.Range("E2").Value() = "'04104013014"

# --- Remove row 6 (the old "suraqa" environment row) entirely, along with its hyperlink ---
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$6') {
        $hl.Delete()
    }
}
$ws.Rows("6:6").Delete()

# --- Update the selected cell shown when the sheet is opened ---
$ws.Range("G10").Select()
